{"js": "// Update each two-digit multiplication answer cell in the table.\n// Old/new text values are each unique within the document, so a direct\n// search-and-replace on the exact string is unambiguous.\nconst replacements = [\n  [\"94\u00d712=1128\", \"38\u00d788=3344\"],\n  [\"88\u00d717=1496\", \"68\u00d725=1700\"],\n  [\"20\u00d738=760\", \"59\u00d754=3186\"],\n  [\"46\u00d760=2760\", \"77\u00d761=4697\"],\n  [\"19\u00d795=1805\", \"83\u00d719=1577\"],\n  [\"53\u00d761=3233\", \"20\u00d795=1900\"],\n  [\"41\u00d794=3854\", \"34\u00d760=2040\"],\n  [\"23\u00d783=1909\", \"63\u00d778=4914\"],\n  [\"49\u00d724=1176\", \"87\u00d736=3132\"],\n  [\"54\u00d718=972\", \"71\u00d740=2840\"],\n  [\"40\u00d736=1440\", \"27\u00d749=1323\"],\n  [\"17\u00d797=1649\", \"16\u00d716=256\"],\n  [\"74\u00d720=1480\", \"66\u00d784=5544\"],\n  [\"74\u00d747=3478\", \"59\u00d712=708\"],\n  [\"56\u00d761=3416\", \"71\u00d793=6603\"],\n  [\"82\u00d791=7462\", \"27\u00d751=1377\"],\n  [\"77\u00d776=5852\", \"59\u00d790=5310\"],\n  [\"92\u00d795=8740\", \"54\u00d740=2160\"],\n  [\"25\u00d752=1300\", \"49\u00d772=3528\"],\n  [\"68\u00d790=6120\", \"94\u00d733=3102\"],\n  [\"67\u00d768=4556\", \"89\u00d759=5251\"],\n  [\"84\u00d726=2184\", \"25\u00d735=875\"],\n  [\"84\u00d739=3276\", \"69\u00d770=4830\"],\n  [\"60\u00d792=5520\", \"42\u00d782=3444\"],\n  [\"58\u00d731=1798\", \"85\u00d794=7990\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "# Update each two-digit multiplication answer cell in the table.\n# Old/new text values are each unique within the document, so Find/Replace\n# on the exact string (wdReplaceAll, but each only matches once) is safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"94\u00d712=1128\", \"38\u00d788=3344\"),\n  @(\"88\u00d717=1496\", \"68\u00d725=1700\"),\n  @(\"20\u00d738=760\", \"59\u00d754=3186\"),\n  @(\"46\u00d760=2760\", \"77\u00d761=4697\"),\n  @(\"19\u00d795=1805\", \"83\u00d719=1577\"),\n  @(\"53\u00d761=3233\", \"20\u00d795=1900\"),\n  @(\"41\u00d794=3854\", \"34\u00d760=2040\"),\n  @(\"23\u00d783=1909\", \"63\u00d778=4914\"),\n  @(\"49\u00d724=1176\", \"87\u00d736=3132\"),\n  @(\"54\u00d718=972\", \"71\u00d740=2840\"),\n  @(\"40\u00d736=1440\", \"27\u00d749=1323\"),\n  @(\"17\u00d797=1649\", \"16\u00d716=256\"),\n  @(\"74\u00d720=1480\", \"66\u00d784=5544\"),\n  @(\"74\u00d747=3478\", \"59\u00d712=708\"),\n  @(\"56\u00d761=3416\", \"71\u00d793=6603\"),\n  @(\"82\u00d791=7462\", \"27\u00d751=1377\"),\n  @(\"77\u00d776=5852\", \"59\u00d790=5310\"),\n  @(\"92\u00d795=8740\", \"54\u00d740=2160\"),\n  @(\"25\u00d752=1300\", \"49\u00d772=3528\"),\n  @(\"68\u00d790=6120\", \"94\u00d733=3102\"),\n  @(\"67\u00d768=4556\", \"89\u00d759=5251\"),\n  @(\"84\u00d726=2184\", \"25\u00d735=875\"),\n  @(\"84\u00d739=3276\", \"69\u00d770=4830\"),\n  @(\"60\u00d792=5520\", \"42\u00d782=3444\"),\n  @(\"58\u00d731=1798\", \"85\u00d794=7990\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $found = $find.Execute(\n    $find.Text,              # FindText\n    $false,                  # MatchCase\n    $false,                  # MatchWholeWord\n    $false,                  # MatchWildcards\n    $false,                  # MatchSoundsLike\n    $false,                  # MatchAllWordForms\n    $true,                   # Forward\n    1,                       # Wrap (wdFindContinue)\n    $false,                  # Format\n    $find.Replacement.Text,  # ReplaceWith\n    2                        # Replace (wdReplaceAll)\n  )\n  if (-not $found) {\n    throw \"Text not found: $oldText\"\n  }\n}"}
